$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 (pushing existing row 26 "003895497" and below down by one)
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new account data.
# The account number has a leading zero, so prefix it with an apostrophe
# (standard Excel text-entry convention) to keep it stored as text instead
# of being auto-converted to a number.
$ws.Cells.Item(26, 1).Value = "'008070544"
$ws.Cells.Item(26, 2).Value = "MARINA"
$ws.Cells.Item(26, 3).Value = 215
